$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate match data among rows 79, 80, 81 (columns F:V), keeping A (index) and E (date) fixed ---
# Capture current (pre-edit) F:V contents for rows 79, 80, 81 using Value2 (Value getter is unreliable here)
$cols = 6..22

$row79 = @{}
$row80 = @{}
$row81 = @{}
foreach ($c in $cols) {
    $row79[$c] = $ws.Cells.Item(79, $c).Value2
    $row80[$c] = $ws.Cells.Item(80, $c).Value2
    $row81[$c] = $ws.Cells.Item(81, $c).Value2
}

# New row79 <- old row81 ; new row80 <- old row79 ; new row81 <- old row80
foreach ($c in $cols) {
    $ws.Cells.Item(79, $c).Value = $row81[$c]
    $ws.Cells.Item(80, $c).Value = $row79[$c]
    $ws.Cells.Item(81, $c).Value = $row80[$c]
}

# --- Append new row 87 ---
# First clone the formatting of row 86 (A:V) onto row 87 so styles (bold/border
# on A, date format on E) match the rest of the table without minting new style ids.
$ws.Range("A86:V86").Copy() | Out-Null
$ws.Range("A87:V87").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "denmark"
$ws.Cells.Item(87, 3).Value = "1st-division"
$ws.Cells.Item(87, 4).Value = "2023-2024"
$ws.Cells.Item(87, 5).Value = 45234.58333333334
$ws.Cells.Item(87, 6).Value = "Vendsyssel"
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = "B.93"
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 1.53
$ws.Cells.Item(87, 11).Value = "29/10/2023 13:13"
$ws.Cells.Item(87, 12).Value = 1.61
$ws.Cells.Item(87, 13).Value = "04/11/2023 13:52"
$ws.Cells.Item(87, 14).Value = 4.36
$ws.Cells.Item(87, 15).Value = "29/10/2023 13:13"
$ws.Cells.Item(87, 16).Value = 4.29
$ws.Cells.Item(87, 17).Value = "04/11/2023 13:52"
$ws.Cells.Item(87, 18).Value = 5.11
$ws.Cells.Item(87, 19).Value = "29/10/2023 13:13"
$ws.Cells.Item(87, 20).Value = 5.12
$ws.Cells.Item(87, 21).Value = "04/11/2023 13:52"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-boldklubben-1893/UZEPq3SH/"
